$d = $word.ActiveDocument

# --- Step 1: insert the new "toastify" paragraph right after the first
#     empty paragraph that follows the "redirected" paragraph (w14:paraId
#     0C1D41C4), i.e. paragraph index 41 in the original document.
$pAfterRedirect = $d.Paragraphs.Item(41)
$rInsert = $pAfterRedirect.Range
$rInsert.Collapse(0)
$rInsert.InsertParagraphAfter()

$pNew = $d.Paragraphs.Item(42)
$xmlNewPara = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">En endring jeg også ville gjort er å legge in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>toastify</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> meldinger som kjøres når ting blir laget, endret og slettet.</w:t></w:r></w:p>'
$pNew.Range.InsertXML($xmlNewPara)

# --- Step 2: insert a brand-new empty paragraph right before the
#     "Det er noen ting ..." paragraph (w14:paraId 30100766), which is now
#     shifted down to index 44 (41 redirect-empty, 42 toastify, 43 second
#     empty 1BD2F63D, 44 "Det er noen ting...").
$pBeforeTarget = $d.Paragraphs.Item(43)
$rInsert2 = $pBeforeTarget.Range
$rInsert2.Collapse(0)
$rInsert2.InsertParagraphAfter()

# --- Step 3: remove the <w:lastRenderedPageBreak/> from the run that
#     begins the "Det er noen ting ..." paragraph (now index 45), while
#     keeping the rest of its run structure (incl. the separate
#     "søkefunksjon" run) untouched.
$pTarget = $d.Paragraphs.Item(45)
$xmlTarget = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="30100766" w14:textId="57B58535" w:rsidR="00F130E4" w:rsidRDefault="00F130E4" w:rsidP="00F130E4"><w:r><w:t xml:space="preserve">Det er noen ting jeg hadde gjort som videreutvikling for oppgaven. En av disse tingene er at jeg hadde lagd en </w:t></w:r><w:r w:rsidR="000C5756"><w:t>søkefunksjon</w:t></w:r><w:r><w:t xml:space="preserve"> for visningen av verktøy. Dette hadde vært fint for hvis bedriften har mye verktøy og du leiter etter en spesifikk type for eksempel så kunne du bare ha søkt det opp i stedet. I samme tema som det så hadde jeg lagd en filtrering som brukeren kunne brukt, dette hadde jeg satt opp så de kunne filtrert på hvor ting er lagret, hva slags type verktøy det er og mellom to datoer for å finne alt som er kjøpt imellom da.</w:t></w:r></w:p>'
$pTarget.Range.InsertXML($xmlTarget)
